$wb = $excel.ActiveWorkbook

# Remove Sheet2 and Sheet3
$wb.Worksheets.Item("Sheet2").Delete()
$wb.Worksheets.Item("Sheet3").Delete()

$ws = $wb.Worksheets.Item("Sheet1")

# Update existing rows with new stat values
$ws.Range("B2").Value = 160
$ws.Range("C2").Value = 180
$ws.Range("D2").Value = 140

$ws.Range("B3").Value = 80
$ws.Range("C3").Value = 320
$ws.Range("D3").Value = 70

# Add new row for BP_EnemySoldier
$ws.Range("A4").Value = "BP_EnemySoldier"
$ws.Range("B4").Value = 90
$ws.Range("C4").Value = 180
$ws.Range("D4").Value = 30
